$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New student rows (3-7) ---
$data = @(
    @("Mark","Bush","Yellow","B","1 Feb, 2022","1 May, 2022","Mon 3:45PM","Tue 5:00PM","Thu 6:00PM"),
    @("Nikki","Bush","Red","A","1 Jan, 2022","1 Apr, 2022","Mon 3:45PM","Wed 6:00PM","N/A"),
    @("Mike","Hawk","Brown","D","1 Jan, 2022","1 Jun, 2022","Tue 6:15PM","Thu 6:00PM","N/A"),
    @("Jill","Brody","Green","D","1 Jan, 2022","1 Jun, 2022","Tue 6:15PM","Wed 6:00PM","Fri 7:30PM"),
    @("Ben","Jamin","No","B","1 Jan, 2022","1 May, 2022","Mon 3:45PM","N/A","N/A")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# --- Widen the class-time columns (G, H, I) ---
$ws.Columns.Item(7).ColumnWidth = 11.833333333333334
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666
$ws.Columns.Item(9).ColumnWidth = 11

# --- Final selection: single cell F6 ---
$ws.Range("F6").Select() | Out-Null
